$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 0.89013549103003475
$ws.Range("AC1").Value = 0.66207240043915927
$ws.Range("AR1").Value = 0.8873312886028093
$ws.Range("AN3").Value = 0.80726888363099336
$ws.Range("AM5").Value = 0.9725121859528959
$ws.Range("AQ5").Value = 0.874992001613762
$ws.Range("BO5").Value = 0.68001571766738511
$ws.Range("Z6").Value = 0.81345373379701247
$ws.Range("AE6").Value = 0.87270707308456319
$ws.Range("I7").Value = 0.79298596809369504
$ws.Range("W7").Value = 0.92207673735410445
$ws.Range("J8").Value = 0.80071516474976956
$ws.Range("AR8").Value = 0.74536481675458954
$ws.Range("O9").Value = 0.63577980489087826
$ws.Range("S9").Value = 0.87477449045373068
$ws.Range("N10").Value = 0.69116486382783959
$ws.Range("AN10").Value = 0.95922751645924076
$ws.Range("M11").Value = 0.73012326293240759
$ws.Range("V11").Value = 0.97749739580224038
$ws.Range("N12").Value = 0.8639012368963066
$ws.Range("AJ12").Value = 0.85907011383312915
$ws.Range("AO12").Value = 0.77982091963826494
$ws.Range("N13").Value = 0.99141467602122535
$ws.Range("O14").Value = 0.698939136859418
$ws.Range("AM14").Value = 0.94020400269072191
$ws.Range("D15").Value = 0.95501951163764787
$ws.Range("AO15").Value = 0.96510653930901569
$ws.Range("AZ15").Value = 0.89548254298502639
$ws.Range("BD15").Value = 0.81524134503086509
$ws.Range("AG16").Value = 0.84707728733676513
$ws.Range("AU16").Value = 0.70424565895922697
$ws.Range("AW17").Value = 0.99907130047479464
$ws.Range("S18").Value = 0.89130154772810877
$ws.Range("Q19").Value = 0.7698552190447927
$ws.Range("BM19").Value = 0.85657208556813702
$ws.Range("G20").Value = 0.75816528074871781
$ws.Range("R20").Value = 0.97038241392971636
$ws.Range("S20").Value = 0.93962495450488115
$ws.Range("W20").Value = 0.66440883216811342
$ws.Range("AG21").Value = 0.95913251337339278
$ws.Range("AY21").Value = 0.83080552707597111
$ws.Range("J22").Value = 0.95152133331046085
$ws.Range("AD22").Value = 0.8934064589303341
$ws.Range("AM22").Value = 0.95799480335733933
$ws.Range("AP22").Value = 0.79557682299644106
$ws.Range("BD22").Value = 0.96000909977905469
$ws.Range("BI22").Value = 0.84440392548514143
$ws.Range("AJ23").Value = 0.85019374560942207
$ws.Range("AL24").Value = 0.79676766481033645
$ws.Range("AS24").Value = 0.90153100348584136
$ws.Range("F25").Value = 0.56954231986337922
$ws.Range("T25").Value = 0.73213758471012258
$ws.Range("V26").Value = 0.78040350107890766
$ws.Range("BK26").Value = 0.68611213971629659
$ws.Range("R27").Value = 0.80454411166495077
$ws.Range("Y27").Value = 0.93450590437568271
$ws.Range("AK27").Value = 0.92522662392404298
$ws.Range("AL27").Value = 0.84067170633062982
$ws.Range("BH28").Value = 0.6270069012099776
$ws.Range("AD29").Value = 0.73507300955594201
$ws.Range("B30").Value = 0.85236262328097068
$ws.Range("AB30").Value = 0.97745787816994811
$ws.Range("AW30").Value = 0.77357583846004807
$ws.Range("AA31").Value = 0.87137755021320284
$ws.Range("AC31").Value = 0.93799904224796427
$ws.Range("BD31").Value = 0.89571766096227878
$ws.Range("BF31").Value = 0.90100974133815259
$ws.Range("P32").Value = 0.64882892655206192
$ws.Range("AM32").Value = 0.98367627497105625
$ws.Range("BF32").Value = 0.85705932285853503
$ws.Range("AD33").Value = 0.72881171473062967
$ws.Range("AT33").Value = 0.66871030155249134
$ws.Range("H34").Value = 0.65762801193672471
$ws.Range("W34").Value = 0.54632915813003302
$ws.Range("X34").Value = 0.99486654975276845
$ws.Range("AJ35").Value = 0.86915705698909163
$ws.Range("BC35").Value = 0.78105383681121354
$ws.Range("X36").Value = 0.97317890224125481
$ws.Range("Z36").Value = 0.86682394065593282
$ws.Range("AT36").Value = 0.87993834713845609
$ws.Range("AL37").Value = 0.72521684719445512
$ws.Range("AU37").Value = 0.61909349029360916
$ws.Range("BA38").Value = 0.99406888606006483
$ws.Range("BB38").Value = 0.58722520763365349
$ws.Range("L39").Value = 0.75014425376374461
$ws.Range("BL39").Value = 0.63418794706579107
$ws.Range("Q40").Value = 0.53579154903736192
$ws.Range("AZ40").Value = 0.68840737541776997
$ws.Range("Q42").Value = 0.90535036652384249
$ws.Range("AN42").Value = 0.89048588758184843
$ws.Range("AW43").Value = 0.57222387731084812
$ws.Range("AN44").Value = 0.79946905799973722
$ws.Range("BP44").Value = 0.84057073923876469
$ws.Range("F45").Value = 0.60238191999491575
$ws.Range("O46").Value = 0.93056368336890316
$ws.Range("AD46").Value = 0.73500559231499474
$ws.Range("BP46").Value = 0.8875028705355017
$ws.Range("B48").Value = 0.99706599992121581
$ws.Range("AO48").Value = 0.95320071873651835
$ws.Range("BG49").Value = 0.95117127154009773
$ws.Range("BN49").Value = 0.99055643862437837
$ws.Range("B50").Value = 0.61461245595806446
$ws.Range("D51").Value = 0.58063047323098438
$ws.Range("AG51").Value = 0.88494151785392794
$ws.Range("BC51").Value = 0.85579221642453074
$ws.Range("AX53").Value = 0.55691819677024057
$ws.Range("AZ53").Value = 0.63043566256223182
$ws.Range("I55").Value = 0.9025855587080579
$ws.Range("AF55").Value = 0.9568922565180491
$ws.Range("AX55").Value = 0.53181574399063047
$ws.Range("BH55").Value = 0.97477183052360084
$ws.Range("BB56").Value = 0.68480943022084861
$ws.Range("AX57").Value = 0.93993196524047007
$ws.Range("M58").Value = 0.89217631380460238
$ws.Range("BE58").Value = 0.91064486024078062
$ws.Range("AP59").Value = 0.91347513687437731
$ws.Range("G61").Value = 0.91635912201162761
$ws.Range("K61").Value = 0.61565839827373114
$ws.Range("BH61").Value = 0.82616725708970673
$ws.Range("O62").Value = 0.68009378507542184
$ws.Range("P62").Value = 0.9150867251196515
$ws.Range("AH62").Value = 0.59282667549271506
$ws.Range("C63").Value = 0.88742534731622902
$ws.Range("AI63").Value = 0.9731698605414334
$ws.Range("BA63").Value = 0.84754890568429708
$ws.Range("BL63").Value = 0.87353130395478396
$ws.Range("V64").Value = 0.74701656089734403
$ws.Range("BJ64").Value = 0.66286405273323168
$ws.Range("AS65").Value = 0.89471522251771907
$ws.Range("AR66").Value = 0.99013952346325129
$ws.Range("BO66").Value = 0.94063051471197912
$ws.Range("AV67").Value = 0.92704040700953061
$ws.Range("BM67").Value = 0.76950470536512849
$ws.Range("BP67").Value = 0.89601012707136629
$ws.Range("P68").Value = 0.96225817613710118
$ws.Range("AA68").Value = 0.78334750575725431
